$wb = $excel.ActiveWorkbook

# Data to append to each sheet: two new rows (159 and 160) continuing the
# daily log pattern already present in the sheet.
#
# Sheet 1: MID_LFT_#1
$ws = $wb.Worksheets.Item(1)
$r = 159
$ws.Cells.Item($r,1).Value = 45945.46252314815
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x01,0x90"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0xC4"
$ws.Cells.Item($r,5).Value = "0x07"
$ws.Cells.Item($r,6).Value = 400
$ws.Cells.Item($r,7).Value = 568631262647113000000000.0
$ws.Cells.Item($r,8).Value = 204
$ws.Cells.Item($r,9).Value = 7

$r = 160
$ws.Cells.Item($r,1).Value = 45946.46368055556
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x01,0x90"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0xC0"
$ws.Cells.Item($r,5).Value = "0x07"
$ws.Cells.Item($r,6).Value = 400
$ws.Cells.Item($r,7).Value = 568631262647113000000000.0
$ws.Cells.Item($r,8).Value = 204
$ws.Cells.Item($r,9).Value = 7

# Sheet 2: MID_LFT_#2
$ws = $wb.Worksheets.Item(2)
$r = 159
$ws.Cells.Item($r,1).Value = 45945.46252314815
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x01,0x7c"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0xE0"
$ws.Cells.Item($r,5).Value = "0x19"
$ws.Cells.Item($r,6).Value = 380
$ws.Cells.Item($r,7).Value = 568432987514711000000000.0
$ws.Cells.Item($r,8).Value = 224
$ws.Cells.Item($r,9).Value = 25

$r = 160
$ws.Cells.Item($r,1).Value = 45946.46368055556
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x01,0x7c"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0xDC"
$ws.Cells.Item($r,5).Value = "0x19"
$ws.Cells.Item($r,6).Value = 380
$ws.Cells.Item($r,7).Value = 568432987514711000000000.0
$ws.Cells.Item($r,8).Value = 220
$ws.Cells.Item($r,9).Value = 25

# Sheet 3: MID_PLT_#1
$ws = $wb.Worksheets.Item(3)
$r = 159
$ws.Cells.Item($r,1).Value = 45945.46252314815
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x00,0x6e"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0x51"
$ws.Cells.Item($r,5).Value = "0x15"
$ws.Cells.Item($r,6).Value = 110
$ws.Cells.Item($r,7).Value = 568631262647113000000000.0
$ws.Cells.Item($r,8).Value = 81
$ws.Cells.Item($r,9).Value = 15

$r = 160
$ws.Cells.Item($r,1).Value = 45946.46368055556
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x00,0x6e"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0x51"
$ws.Cells.Item($r,5).Value = "0x15"
$ws.Cells.Item($r,6).Value = 110
$ws.Cells.Item($r,7).Value = 568631262647113000000000.0
$ws.Cells.Item($r,8).Value = 81
$ws.Cells.Item($r,9).Value = 15

# Sheet 4: MID_PLT_#2
$ws = $wb.Worksheets.Item(4)
$r = 159
$ws.Cells.Item($r,1).Value = 45945.46252314815
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x00,0x82"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0x6A"
$ws.Cells.Item($r,5).Value = "0x9"
$ws.Cells.Item($r,6).Value = 130
$ws.Cells.Item($r,7).Value = 568631262647113000000000.0
$ws.Cells.Item($r,8).Value = 106
$ws.Cells.Item($r,9).Value = 9

$r = 160
$ws.Cells.Item($r,1).Value = 45946.46368055556
$ws.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,2).Value = "0x00,0x82"
$ws.Cells.Item($r,3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item($r,4).Value = "0x00,0x69"
$ws.Cells.Item($r,5).Value = "0x9"
$ws.Cells.Item($r,6).Value = 130
$ws.Cells.Item($r,7).Value = 568631262647113000000000.0
$ws.Cells.Item($r,8).Value = 105
$ws.Cells.Item($r,9).Value = 9
